$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.731.75'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.477.90'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.20'
$ws.Range('E5').Value = '  +1.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.17'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.511'
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0859'
$ws.Range('E10').Value = '  +2.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '33.04'
$ws.Range('E11').Value = '  +1.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.109'
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.860.71'
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.89'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.51'
$ws.Range('E15').Value = '  -1.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.482.91'
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.794'
$ws.Range('E17').Value = '  +2.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.696.15'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0942'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.08'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.26'
$ws.Range('E22').Value = '  -1.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.79'
$ws.Range('E23').Value = '  +0.52%  '
$ws.Range('E24').Value = '  +1.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.94'
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.98'
$ws.Range('E27').Value = '  +2.13%  '
$ws.Range('E28').Value = '  -1.15%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.63'
$ws.Range('E30').Value = '  +4.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.34'
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.44'
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0765'
$ws.Range('E34').Value = '  +0.66%  '
$ws.Range('E35').Value = '  -0.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.16'
$ws.Range('E36').Value = '  -1.49%  '
$ws.Range('E37').Value = '  +2.97%  '
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('E39').Value = '  -0.37%  '
$ws.Range('E40').Value = '  +0.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.00'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.43'
$ws.Range('E42').Value = '  -1.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.997.15'
$ws.Range('E43').Value = '  +1.25%  '
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.66'
$ws.Range('E45').Value = '  -1.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.97'
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.44'
$ws.Range('E47').Value = '  +4.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.739.11'
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('E49').Value = '  +5.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '97.70'
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '67.31'
$ws.Range('E51').Value = '  +0.49%  '
